# Add two new list paragraphs at the end of the document body, right
# after the last existing paragraph ("Các báo cáo chi tiết ...") and
# before the section break:
#   1. A bold heading paragraph (ListParagraph style, numId=1):
#        "Quản lý dự án và dịch vụ (Project and services management)"
#   2. A body paragraph (ListParagraph style, numId=2) whose paragraph
#      mark keeps bold formatting but whose run text is not bold:
#        "Các dịch vụ có liên quan đến hóa đơn, cần nhân lực, tài nguyên."

$d = $word.ActiveDocument

# Insertion point: the very end of the document body (just before sectPr).
$endPos = $d.Content.End
$insertionRange = $d.Range($endPos, $endPos)

$newParagraphsXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Quản lý dự án và dịch vụ (Project and services management)</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:t>Các dịch vụ có liên quan đến hóa đơn, cần nhân lực, tài nguyên.</w:t></w:r></w:p>'

$insertionRange.InsertXML($newParagraphsXml)
